$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22 for the new "BitDAO" entry,
# shifting existing rows 22-51 down to 23-52 (mirrors source list update).
$ws.Rows.Item(22).Insert()

# Remove the now-duplicated last row (52) so the table stays A1:E51,
# dropping the previously-last entry (EnergySwap) which fell off the list.
$ws.Rows.Item(52).Delete()

# Force columns B:E to Text format so numeric-looking strings (prices,
# percentages) are stored as text, matching the source data (inlineStr).
$ws.Range('B2:E51').NumberFormat = '@'

# Re-assert every data cell (A2:E51) to match the refreshed crypto listing.
$ws.Range('A2').Value() = 0
$ws.Range('B2').Value() = 'Bitcoin'
$ws.Range('C2').Value() = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range('D2').Value() = '30.491.28'
$ws.Range('E2').Value() = '  -1.14%  '

$ws.Range('A3').Value() = 1
$ws.Range('B3').Value() = 'Ethereum'
$ws.Range('C3').Value() = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range('D3').Value() = '1.912.55'
$ws.Range('E3').Value() = '  -1.30%  '

$ws.Range('A4').Value() = 2
$ws.Range('B4').Value() = 'TetherUSD'
$ws.Range('C4').Value() = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range('D4').Value() = '1.001'
$ws.Range('E4').Value() = '  +0.10%  '

$ws.Range('A5').Value() = 3
$ws.Range('B5').Value() = 'BNB'
$ws.Range('C5').Value() = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').Value() = '239.91'
$ws.Range('E5').Value() = '  -1.47%  '

$ws.Range('A6').Value() = 4
$ws.Range('B6').Value() = 'USDC'
$ws.Range('C6').Value() = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D6').Value() = '1.001'
$ws.Range('E6').Value() = '  +0.10%  '

$ws.Range('A7').Value() = 5
$ws.Range('B7').Value() = 'XRP'
$ws.Range('C7').Value() = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D7').Value() = '0.4786'
$ws.Range('E7').Value() = '  -2.51%  '

$ws.Range('A8').Value() = 6
$ws.Range('B8').Value() = 'Cardano'
$ws.Range('C8').Value() = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').Value() = '0.2844'
$ws.Range('E8').Value() = '  -3.10%  '

$ws.Range('A9').Value() = 7
$ws.Range('B9').Value() = 'Dogecoin'
$ws.Range('C9').Value() = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').Value() = '0.06698'
$ws.Range('E9').Value() = '  -2.64%  '

$ws.Range('A10').Value() = 8
$ws.Range('B10').Value() = 'Solana'
$ws.Range('C10').Value() = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D10').Value() = '19.41'
$ws.Range('E10').Value() = '  +0.93%  '

$ws.Range('A11').Value() = 9
$ws.Range('B11').Value() = 'Litecoin'
$ws.Range('C11').Value() = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D11').Value() = '102.88'
$ws.Range('E11').Value() = '  -1.93%  '

$ws.Range('A12').Value() = 10
$ws.Range('B12').Value() = 'TRON'
$ws.Range('C12').Value() = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Value() = '0.07794'
$ws.Range('E12').Value() = '  +0.27%  '

$ws.Range('A13').Value() = 11
$ws.Range('B13').Value() = 'WrappedEther'
$ws.Range('C13').Value() = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value() = '1.921.28'
$ws.Range('E13').Value() = '  -0.53%  '

$ws.Range('A14').Value() = 12
$ws.Range('B14').Value() = 'Polkadot'
$ws.Range('C14').Value() = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value() = '5.202'
$ws.Range('E14').Value() = '  -2.81%  '

$ws.Range('A15').Value() = 13
$ws.Range('B15').Value() = 'Polygon'
$ws.Range('C15').Value() = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').Value() = '0.6702'
$ws.Range('E15').Value() = '  -4.06%  '

$ws.Range('A16').Value() = 14
$ws.Range('B16').Value() = 'BitcoinCash'
$ws.Range('C16').Value() = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D16').Value() = '276.35'
$ws.Range('E16').Value() = '  +0.85%  '

$ws.Range('A17').Value() = 15
$ws.Range('B17').Value() = 'WrappedBTC'
$ws.Range('C17').Value() = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value() = '30.539.35'
$ws.Range('E17').Value() = '  -1.03%  '

$ws.Range('A18').Value() = 16
$ws.Range('B18').Value() = 'Dai'
$ws.Range('C18').Value() = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D18').Value() = '1.001'
$ws.Range('E18').Value() = '  -0.02%  '

$ws.Range('A19').Value() = 17
$ws.Range('B19').Value() = 'ShibaInu'
$ws.Range('C19').Value() = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').Value() = '0.000007481'
$ws.Range('E19').Value() = '  -3.09%  '

$ws.Range('A20').Value() = 18
$ws.Range('B20').Value() = 'Avalanche'
$ws.Range('C20').Value() = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').Value() = '12.63'
$ws.Range('E20').Value() = '  -3.44%  '

$ws.Range('A21').Value() = 19
$ws.Range('B21').Value() = 'Uniswap'
$ws.Range('C21').Value() = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value() = '5.387'
$ws.Range('E21').Value() = '  -3.35%  '

$ws.Range('A22').Value() = 20
$ws.Range('B22').Value() = 'BitDAO'
$ws.Range('C22').Value() = 'https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit'
$ws.Range('D22').Value() = '0.4691'
$ws.Range('E22').Value() = '  -6.24%  '

$ws.Range('A23').Value() = 21
$ws.Range('B23').Value() = 'BinanceUSD'
$ws.Range('C23').Value() = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D23').Value() = '1.001'
$ws.Range('E23').Value() = '  +0.02%  '

$ws.Range('A24').Value() = 22
$ws.Range('B24').Value() = 'Chainlink'
$ws.Range('C24').Value() = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D24').Value() = '6.302'
$ws.Range('E24').Value() = '  -3.44%  '

$ws.Range('A25').Value() = 23
$ws.Range('B25').Value() = 'Cosmos'
$ws.Range('C25').Value() = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D25').Value() = '9.340'
$ws.Range('E25').Value() = '  -5.26%  '

$ws.Range('A26').Value() = 24
$ws.Range('B26').Value() = 'Monero'
$ws.Range('C26').Value() = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value() = '167.52'
$ws.Range('E26').Value() = '  +0.88%  '

$ws.Range('A27').Value() = 25
$ws.Range('B27').Value() = 'EthereumClassic'
$ws.Range('C27').Value() = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value() = '19.19'
$ws.Range('E27').Value() = '  -2.09%  '

$ws.Range('A28').Value() = 26
$ws.Range('B28').Value() = 'LidoDAOToken'
$ws.Range('C28').Value() = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').Value() = '2.080'
$ws.Range('E28').Value() = '  -3.61%  '

$ws.Range('A29').Value() = 27
$ws.Range('B29').Value() = 'Toncoin'
$ws.Range('C29').Value() = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').Value() = '1.384'
$ws.Range('E29').Value() = '  -0.56%  '

$ws.Range('A30').Value() = 28
$ws.Range('B30').Value() = 'Stellar'
$ws.Range('C30').Value() = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').Value() = '0.09953'
$ws.Range('E30').Value() = '  -4.53%  '

$ws.Range('A31').Value() = 29
$ws.Range('B31').Value() = 'Filecoin'
$ws.Range('C31').Value() = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value() = '4.584'
$ws.Range('E31').Value() = '  +0.72%  '

$ws.Range('A32').Value() = 30
$ws.Range('B32').Value() = 'PancakeSwap'
$ws.Range('C32').Value() = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value() = '1.517'
$ws.Range('E32').Value() = '  -2.52%  '

$ws.Range('A33').Value() = 31
$ws.Range('B33').Value() = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value() = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').Value() = '4.251'
$ws.Range('E33').Value() = '  -2.89%  '

$ws.Range('A34').Value() = 32
$ws.Range('B34').Value() = 'Hedera'
$ws.Range('C34').Value() = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value() = '0.04733'
$ws.Range('E34').Value() = '  -3.11%  '

$ws.Range('A35').Value() = 33
$ws.Range('B35').Value() = 'ImmutableX'
$ws.Range('C35').Value() = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value() = '0.7249'
$ws.Range('E35').Value() = '  -4.35%  '

$ws.Range('A36').Value() = 34
$ws.Range('B36').Value() = 'ARBITRUM'
$ws.Range('C36').Value() = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value() = '1.111'
$ws.Range('E36').Value() = '  -3.35%  '

$ws.Range('A37').Value() = 35
$ws.Range('B37').Value() = 'HuobiToken'
$ws.Range('C37').Value() = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').Value() = '2.716'
$ws.Range('E37').Value() = '  -0.71%  '

$ws.Range('A38').Value() = 36
$ws.Range('B38').Value() = 'VeChain'
$ws.Range('C38').Value() = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value() = '0.01906'
$ws.Range('E38').Value() = '  -5.04%  '

$ws.Range('A39').Value() = 37
$ws.Range('B39').Value() = 'MXToken'
$ws.Range('C39').Value() = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value() = '2.623'
$ws.Range('E39').Value() = '  -1.24%  '

$ws.Range('A40').Value() = 38
$ws.Range('B40').Value() = 'FraxShare'
$ws.Range('C40').Value() = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value() = '6.321'
$ws.Range('E40').Value() = '  -3.38%  '

$ws.Range('A41').Value() = 39
$ws.Range('B41').Value() = 'Aave'
$ws.Range('C41').Value() = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value() = '73.82'
$ws.Range('E41').Value() = '  -5.01%  '

$ws.Range('A42').Value() = 40
$ws.Range('B42').Value() = 'RenderToken'
$ws.Range('C42').Value() = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').Value() = '1.954'
$ws.Range('E42').Value() = '  -6.53%  '

$ws.Range('A43').Value() = 41
$ws.Range('B43').Value() = 'TrustWalletToken'
$ws.Range('C43').Value() = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value() = '0.8618'
$ws.Range('E43').Value() = '  -5.53%  '

$ws.Range('A44').Value() = 42
$ws.Range('B44').Value() = 'Quant'
$ws.Range('C44').Value() = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').Value() = '106.16'
$ws.Range('E44').Value() = '  -1.52%  '

$ws.Range('A45').Value() = 43
$ws.Range('B45').Value() = 'TheSandbox'
$ws.Range('C45').Value() = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').Value() = '0.4256'
$ws.Range('E45').Value() = '  -4.05%  '

$ws.Range('A46').Value() = 44
$ws.Range('B46').Value() = 'PaxDollar'
$ws.Range('C46').Value() = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D46').Value() = '1.001'
$ws.Range('E46').Value() = '  +0.23%  '

$ws.Range('A47').Value() = 45
$ws.Range('B47').Value() = 'Aptos'
$ws.Range('C47').Value() = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D47').Value() = '7.389'
$ws.Range('E47').Value() = '  -3.94%  '

$ws.Range('A48').Value() = 46
$ws.Range('B48').Value() = 'Maker'
$ws.Range('C48').Value() = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value() = '955.67'
$ws.Range('E48').Value() = '  -3.94%  '

$ws.Range('A49').Value() = 47
$ws.Range('B49').Value() = 'Algorand'
$ws.Range('C49').Value() = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value() = '0.1202'
$ws.Range('E49').Value() = '  -3.79%  '

$ws.Range('A50').Value() = 48
$ws.Range('B50').Value() = 'Elrond'
$ws.Range('C50').Value() = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').Value() = '34.63'
$ws.Range('E50').Value() = '  -4.06%  '

$ws.Range('A51').Value() = 49
$ws.Range('B51').Value() = 'Cronos'
$ws.Range('C51').Value() = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value() = '0.05801'
$ws.Range('E51').Value() = '  +0.59%  '
